$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 261
$ws.Range("F4").Value = 589
$ws.Range("F5").Value = 2582
$ws.Range("F6").Value = 11
$ws.Range("F7").Value = 173
$ws.Range("F9").Value = 252
$ws.Range("F10").Value = 5413
$ws.Range("F11").Value = 103
$ws.Range("F12").Value = 1473
$ws.Range("F13").Value = 1391
$ws.Range("F14").Value = 599
$ws.Range("F15").Value = 6991
$ws.Range("F17").Value = 50
$ws.Range("F20").Value = 4688
$ws.Range("F22").Value = 2350
$ws.Range("F23").Value = 1259
$ws.Range("F24").Value = 445
$ws.Range("F25").Value = 1157
$ws.Range("F26").Value = 219
$ws.Range("F27").Value = 91
$ws.Range("F28").Value = 74
$ws.Range("F29").Value = 165
$ws.Range("F30").Value = 367
$ws.Range("F31").Value = 1282
$ws.Range("F32").Value = 1992
$ws.Range("F33").Value = 240
$ws.Range("F34").Value = 525
$ws.Range("F35").Value = 2
$ws.Range("F36").Value = 205
$ws.Range("F37").Value = 1373
$ws.Range("F41").Value = 175
$ws.Range("F42").Value = 1124
$ws.Range("F43").Value = 2411
$ws.Range("F44").Value = 42
$ws.Range("F45").Value = 69
$ws.Range("F47").Value = 234
$ws.Range("F49").Value = 19

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 455
$ws.Range("F6").Value = 246
$ws.Range("F7").Value = 132
$ws.Range("F10").Value = 4
$ws.Range("F16").Value = 185
$ws.Range("G16").Value = 304
$ws.Range("F23").Value = 140
$ws.Range("F24").Value = 32
$ws.Range("F28").Value = 286

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 496
$ws.Range("F6").Value = 1657
$ws.Range("F7").Value = 537
$ws.Range("F8").Value = 1302
$ws.Range("F10").Value = 1745
$ws.Range("F11").Value = 2200
$ws.Range("F12").Value = 639
$ws.Range("F13").Value = 532

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 261
$ws.Range("F3").Value = 537
$ws.Range("F4").Value = 2582
$ws.Range("F5").Value = 173
$ws.Range("F6").Value = 1302
$ws.Range("F7").Value = 252
$ws.Range("F8").Value = 2200
$ws.Range("F9").Value = 5413
$ws.Range("F10").Value = 639
$ws.Range("F11").Value = 455
$ws.Range("F12").Value = 246
$ws.Range("F13").Value = 132
$ws.Range("F14").Value = 103
$ws.Range("F16").Value = 1473
$ws.Range("F17").Value = 1391
$ws.Range("F18").Value = 6991
$ws.Range("F21").Value = 4
$ws.Range("F22").Value = 50
$ws.Range("F23").Value = 4688
$ws.Range("F24").Value = 2350
$ws.Range("F25").Value = 1157
$ws.Range("F26").Value = 91
$ws.Range("F27").Value = 74
$ws.Range("F30").Value = 165
$ws.Range("F32").Value = 185
$ws.Range("G32").Value = 304
$ws.Range("F33").Value = 367
$ws.Range("F34").Value = 1282
$ws.Range("F35").Value = 1992
$ws.Range("F36").Value = 525
$ws.Range("F38").Value = 205
$ws.Range("F39").Value = 1373
$ws.Range("F44").Value = 1124
$ws.Range("F45").Value = 2411
$ws.Range("F46").Value = 69
$ws.Range("F47").Value = 234
$ws.Range("F49").Value = 19
